# Scheduled market-price refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit figures across all job sheets to match the
# latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2500
$ws.Range("I18").Value = 1428.5714
$ws.Range("K18").Value = 1428.5714
$ws.Range("M18").Value = -1144.5714

$ws.Range("H86").Value = 3713.5715
$ws.Range("I86").Value = 2650
$ws.Range("K86").Value = 2650
$ws.Range("M86").Value = -1527

$ws.Range("H89").Value = 3713.5715
$ws.Range("I89").Value = 2650
$ws.Range("K89").Value = 13250
$ws.Range("M89").Value = -7634

$ws.Range("H137").Value = 1464.2646
$ws.Range("I137").Value = 1112.1111
$ws.Range("K137").Value = 3336.3333
$ws.Range("M137").Value = -786.3333000000002

$ws.Range("H140").Value = 96249.5
$ws.Range("J140").Value = 96249.5
$ws.Range("L140").Value = 96249.5
$ws.Range("N140").Value = -106609.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 888.7778
$ws.Range("I2").Value = 888.7778
$ws.Range("K2").Value = 888.7778
$ws.Range("M2").Value = -775.7778

$ws.Range("H45").Value = 5317.48
$ws.Range("I45").Value = 5795
$ws.Range("K45").Value = 5795
$ws.Range("M45").Value = -5418

$ws.Range("H61").Value = 4666.436
$ws.Range("I61").Value = 3499.4285
$ws.Range("K61").Value = 3499.4285
$ws.Range("M61").Value = -3287.4285

$ws.Range("H74").Value = 4333.657
$ws.Range("I74").Value = 1117.3182
$ws.Range("J74").Value = 9776.691999999999
$ws.Range("K74").Value = 1117.3182
$ws.Range("L74").Value = 9776.691999999999
$ws.Range("M74").Value = -243.3181999999999
$ws.Range("N74").Value = -11524.692

$ws.Range("H77").Value = 4333.657
$ws.Range("I77").Value = 1117.3182
$ws.Range("J77").Value = 9776.691999999999
$ws.Range("K77").Value = 5586.590999999999
$ws.Range("L77").Value = 48883.45999999999
$ws.Range("M77").Value = -1218.590999999999
$ws.Range("N77").Value = -57619.45999999999

$ws.Range("H116").Value = 888.7778
$ws.Range("I116").Value = 888.7778
$ws.Range("K116").Value = 888.7778
$ws.Range("M116").Value = 1405.2222

$ws.Range("H136").Value = 4666.436
$ws.Range("I136").Value = 3499.4285
$ws.Range("K136").Value = 10498.2855
$ws.Range("M136").Value = -7948.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 888.7778
$ws.Range("I3").Value = 888.7778
$ws.Range("K3").Value = 888.7778
$ws.Range("M3").Value = -774.7778

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H134").Value = 1960.9744
$ws.Range("I134").Value = 1937.9296
$ws.Range("K134").Value = 5813.7888
$ws.Range("M134").Value = -3278.7888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 380

$ws.Range("H58").Value = 4424.32
$ws.Range("I58").Value = 4054.353
$ws.Range("K58").Value = 4054.353
$ws.Range("M58").Value = -3851.353

$ws.Range("H134").Value = 24906.818
$ws.Range("I134").Value = 14332.223
$ws.Range("J134").Value = 72492.5
$ws.Range("K134").Value = 42996.669
$ws.Range("L134").Value = 217477.5
$ws.Range("M134").Value = -40461.669
$ws.Range("N134").Value = -222547.5

$ws.Range("H136").Value = 4424.32
$ws.Range("I136").Value = 4054.353
$ws.Range("K136").Value = 12163.059
$ws.Range("M136").Value = -9613.059000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 541.2
$ws.Range("J11").Value = 176.5
$ws.Range("L11").Value = 529.5
$ws.Range("N11").Value = -809.5

$ws.Range("H29").Value = 1509.8334
$ws.Range("J29").Value = 2866.3333
$ws.Range("L29").Value = 8598.999899999999
$ws.Range("N29").Value = -9152.999899999999

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H140").Value = 2190.9092
$ws.Range("I140").Value = 1410
$ws.Range("J140").Value = 10000
$ws.Range("K140").Value = 4230
$ws.Range("L140").Value = 30000
$ws.Range("M140").Value = 950
$ws.Range("N140").Value = -40360

$ws.Range("H141").Value = 77788.07000000001
$ws.Range("J141").Value = 83387.16
$ws.Range("L141").Value = 250161.48
$ws.Range("N141").Value = -260521.48

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I52").Value = 60000
$ws.Range("K52").Value = 60000
$ws.Range("M52").Value = -59741

$ws.Range("H97").Value = 1084.6364
$ws.Range("I97").Value = 810.375
$ws.Range("K97").Value = 810.375
$ws.Range("M97").Value = -314.375

$ws.Range("H102").Value = 50000864
$ws.Range("I102").Value = 925.625
$ws.Range("K102").Value = 925.625
$ws.Range("M102").Value = 696.375

$ws.Range("H122").Value = 6182
$ws.Range("I122").Value = 5779.2
$ws.Range("K122").Value = 17337.6
$ws.Range("M122").Value = -14887.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3721.4707
$ws.Range("I7").Value = 2183
$ws.Range("K7").Value = 2183
$ws.Range("M7").Value = -2071

$ws.Range("H22").Value = 3285.875
$ws.Range("J22").Value = 3697.8
$ws.Range("L22").Value = 3697.8
$ws.Range("N22").Value = -4287.8

$ws.Range("H27").Value = 3285.875
$ws.Range("J27").Value = 3697.8
$ws.Range("L27").Value = 3697.8
$ws.Range("N27").Value = -3911.8

$ws.Range("H40").Value = 2412.0513
$ws.Range("I40").Value = 1623.0646
$ws.Range("K40").Value = 1623.0646
$ws.Range("M40").Value = -1487.0646

$ws.Range("H61").Value = 2651.3333
$ws.Range("I61").Value = 2444.6428
$ws.Range("J61").Value = 3374.75
$ws.Range("K61").Value = 2444.6428
$ws.Range("L61").Value = 3374.75
$ws.Range("M61").Value = -2242.6428
$ws.Range("N61").Value = -3778.75

$ws.Range("H68").Value = 3031
$ws.Range("I68").Value = 2999.7144
$ws.Range("K68").Value = 2999.7144
$ws.Range("M68").Value = -2250.7144

$ws.Range("H71").Value = 3031
$ws.Range("I71").Value = 2999.7144
$ws.Range("K71").Value = 14998.572
$ws.Range("M71").Value = -11254.572

$ws.Range("H113").Value = 2651.3333
$ws.Range("I113").Value = 2444.6428
$ws.Range("J113").Value = 3374.75
$ws.Range("K113").Value = 2444.6428
$ws.Range("L113").Value = 3374.75
$ws.Range("M113").Value = -274.6428000000001
$ws.Range("N113").Value = -7714.75

$ws.Range("H122").Value = 5039.4116
$ws.Range("I122").Value = 4233.8184
$ws.Range("J122").Value = 6516.3335
$ws.Range("K122").Value = 12701.4552
$ws.Range("L122").Value = 19549.0005
$ws.Range("M122").Value = -10251.4552
$ws.Range("N122").Value = -24449.0005

$ws.Range("H126").Value = 3721.4707
$ws.Range("I126").Value = 2183
$ws.Range("K126").Value = 6549
$ws.Range("M126").Value = -4079

$ws.Range("H132").Value = 4046.1738
$ws.Range("I132").Value = 4097.9375
$ws.Range("J132").Value = 3927.8572
$ws.Range("K132").Value = 12293.8125
$ws.Range("L132").Value = 11783.5716
$ws.Range("M132").Value = -9763.8125
$ws.Range("N132").Value = -16843.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7513166.5
$ws.Range("I15").Value = 9012000
$ws.Range("K15").Value = 9012000
$ws.Range("M15").Value = -9011712

$ws.Range("H62").Value = 7954.3335
$ws.Range("I62").Value = 1509
$ws.Range("J62").Value = 14399.667
$ws.Range("K62").Value = 1509
$ws.Range("L62").Value = 14399.667
$ws.Range("M62").Value = -885
$ws.Range("N62").Value = -15647.667

$ws.Range("H65").Value = 7954.3335
$ws.Range("I65").Value = 1509
$ws.Range("J65").Value = 14399.667
$ws.Range("K65").Value = 7545
$ws.Range("L65").Value = 71998.33499999999
$ws.Range("M65").Value = -4425
$ws.Range("N65").Value = -78238.33499999999

$ws.Range("H122").Value = 1503.0857
$ws.Range("I122").Value = 1409.9333
$ws.Range("K122").Value = 4229.7999
$ws.Range("M122").Value = -1779.7999

$ws.Range("H126").Value = 2527.818
$ws.Range("I126").Value = 2376
$ws.Range("K126").Value = 7128
$ws.Range("M126").Value = -4658

$ws.Range("H132").Value = 2580.0588
$ws.Range("I132").Value = 2580.0588
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7740.176399999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5210.176399999999
$ws.Range("N132").ClearContents()

$ws.Range("I136").Value = 63993.188
$ws.Range("K136").Value = 191979.564
$ws.Range("M136").Value = -189429.564
